$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.911561666666666
$ws.Cells.Item(2, 8).Value = 5.734684999999999
$ws.Cells.Item(2, 9).Value = 0.1720155802183755
$ws.Cells.Item(2, 10).Value = 0.1720155802183755
$ws.Cells.Item(2, 13).Value = 14.25737566666667
$ws.Cells.Item(2, 14).Value = 42.772127
$ws.Cells.Item(2, 15).Value = 0.2087950866344732
$ws.Cells.Item(2, 16).Value = 0.2087950866344732
$ws.Cells.Item(2, 17).Value = 27.25385279166611
$ws.Cells.Item(2, 18).Value = 245.284675124995
$ws.Cells.Item(2, 19).Value = 0.03591600797417489
$ws.Cells.Item(2, 20).Value = 0.03591600797417488
$ws.Cells.Item(3, 7).Value = 1.911561666666666
$ws.Cells.Item(3, 8).Value = 5.734684999999999
$ws.Cells.Item(3, 9).Value = 0.1720155802183755
$ws.Cells.Item(3, 10).Value = 0.1720155802183755
$ws.Cells.Item(3, 14).Value = 87.128332
$ws.Cells.Item(3, 15).Value = 0.4253229592313036
$ws.Cells.Item(3, 16).Value = 0.4253229592313036
$ws.Cells.Item(3, 17).Value = 55.51705984393555
$ws.Cells.Item(3, 18).Value = 499.6535385954199
$ws.Cells.Item(3, 19).Value = 0.07316217561236917
$ws.Cells.Item(3, 20).Value = 0.07316217561236915
$ws.Cells.Item(4, 7).Value = 1.911561666666666
$ws.Cells.Item(4, 8).Value = 5.734684999999999
$ws.Cells.Item(4, 9).Value = 0.1720155802183755
$ws.Cells.Item(4, 10).Value = 0.1720155802183755
$ws.Cells.Item(4, 13).Value = 20.11084633333333
$ws.Cells.Item(4, 14).Value = 60.332539
$ws.Cells.Item(4, 15).Value = 0.2945174484164121
$ws.Cells.Item(4, 16).Value = 0.2945174484164122
$ws.Cells.Item(4, 17).Value = 38.44312293502388
$ws.Cells.Item(4, 18).Value = 345.9881064152149
$ws.Cells.Item(4, 19).Value = 0.05066158977378461
$ws.Cells.Item(4, 20).Value = 0.05066158977378461
$ws.Cells.Item(5, 7).Value = 1.911561666666666
$ws.Cells.Item(5, 8).Value = 5.734684999999999
$ws.Cells.Item(5, 9).Value = 0.1720155802183755
$ws.Cells.Item(5, 10).Value = 0.1720155802183755
$ws.Cells.Item(5, 13).Value = 4.873057999999999
$ws.Cells.Item(5, 14).Value = 14.619174
$ws.Cells.Item(5, 15).Value = 0.07136450571781097
$ws.Cells.Item(5, 16).Value = 0.07136450571781099
$ws.Cells.Item(5, 17).Value = 9.315150872243331
$ws.Cells.Item(5, 18).Value = 83.83635785018998
$ws.Cells.Item(5, 19).Value = 0.01227580685804683
$ws.Cells.Item(5, 20).Value = 0.01227580685804683
$ws.Cells.Item(6, 9).Value = 0.4009917520372743
$ws.Cells.Item(6, 10).Value = 0.4009917520372743
$ws.Cells.Item(6, 13).Value = 14.25737566666667
$ws.Cells.Item(6, 14).Value = 42.772127
$ws.Cells.Item(6, 15).Value = 0.2087950866344732
$ws.Cells.Item(6, 16).Value = 0.2087950866344732
$ws.Cells.Item(6, 17).Value = 63.53244378690712
$ws.Cells.Item(6, 18).Value = 571.7919940821641
$ws.Cells.Item(6, 19).Value = 0.08372510760633188
$ws.Cells.Item(6, 20).Value = 0.08372510760633188
$ws.Cells.Item(7, 9).Value = 0.4009917520372743
$ws.Cells.Item(7, 10).Value = 0.4009917520372743
$ws.Cells.Item(7, 14).Value = 87.128332
$ws.Cells.Item(7, 15).Value = 0.4253229592313036
$ws.Cells.Item(7, 16).Value = 0.4253229592313036
$ws.Cells.Item(7, 17).Value = 129.4178298646916
$ws.Cells.Item(7, 19).Value = 0.1705509986038386
$ws.Cells.Item(7, 20).Value = 0.1705509986038386
$ws.Cells.Item(8, 9).Value = 0.4009917520372743
$ws.Cells.Item(8, 10).Value = 0.4009917520372743
$ws.Cells.Item(8, 13).Value = 20.11084633333333
$ws.Cells.Item(8, 14).Value = 60.332539
$ws.Cells.Item(8, 15).Value = 0.2945174484164121
$ws.Cells.Item(8, 16).Value = 0.2945174484164122
$ws.Cells.Item(8, 17).Value = 89.61615686166088
$ws.Cells.Item(8, 18).Value = 806.545411754948
$ws.Cells.Item(8, 19).Value = 0.1180990676460447
$ws.Cells.Item(8, 20).Value = 0.1180990676460447
$ws.Cells.Item(9, 9).Value = 0.4009917520372743
$ws.Cells.Item(9, 10).Value = 0.4009917520372743
$ws.Cells.Item(9, 13).Value = 4.873057999999999
$ws.Cells.Item(9, 14).Value = 14.619174
$ws.Cells.Item(9, 15).Value = 0.07136450571781097
$ws.Cells.Item(9, 16).Value = 0.07136450571781099
$ws.Cells.Item(9, 17).Value = 21.71488573308533
$ws.Cells.Item(9, 18).Value = 195.433971597768
$ws.Cells.Item(9, 19).Value = 0.0286165781810591
$ws.Cells.Item(9, 20).Value = 0.02861657818105911
$ws.Cells.Item(10, 7).Value = 4.603447666666667
$ws.Cells.Item(10, 8).Value = 13.810343
$ws.Cells.Item(10, 9).Value = 0.4142501574471451
$ws.Cells.Item(10, 10).Value = 0.4142501574471449
$ws.Cells.Item(10, 13).Value = 14.25737566666667
$ws.Cells.Item(10, 14).Value = 42.772127
$ws.Cells.Item(10, 15).Value = 0.2087950866344732
$ws.Cells.Item(10, 16).Value = 0.2087950866344732
$ws.Cells.Item(10, 17).Value = 65.63308274550678
$ws.Cells.Item(10, 18).Value = 590.697744709561
$ws.Cells.Item(10, 19).Value = 0.0864933975125208
$ws.Cells.Item(10, 20).Value = 0.0864933975125208
$ws.Cells.Item(11, 7).Value = 4.603447666666667
$ws.Cells.Item(11, 8).Value = 13.810343
$ws.Cells.Item(11, 9).Value = 0.4142501574471451
$ws.Cells.Item(11, 10).Value = 0.4142501574471449
$ws.Cells.Item(11, 14).Value = 87.128332
$ws.Cells.Item(11, 15).Value = 0.4253229592313036
$ws.Cells.Item(11, 16).Value = 0.4253229592313036
$ws.Cells.Item(11, 17).Value = 133.6969055486529
$ws.Cells.Item(11, 18).Value = 1203.272149937876
$ws.Cells.Item(11, 19).Value = 0.1761901028274532
$ws.Cells.Item(11, 20).Value = 0.1761901028274531
$ws.Cells.Item(12, 7).Value = 4.603447666666667
$ws.Cells.Item(12, 8).Value = 13.810343
$ws.Cells.Item(12, 9).Value = 0.4142501574471451
$ws.Cells.Item(12, 10).Value = 0.4142501574471449
$ws.Cells.Item(12, 13).Value = 20.11084633333333
$ws.Cells.Item(12, 14).Value = 60.332539
$ws.Cells.Item(12, 15).Value = 0.2945174484164121
$ws.Cells.Item(12, 16).Value = 0.2945174484164122
$ws.Cells.Item(12, 17).Value = 92.57922862787521
$ws.Cells.Item(12, 18).Value = 833.213057650877
$ws.Cells.Item(12, 19).Value = 0.1220038993774302
$ws.Cells.Item(12, 20).Value = 0.1220038993774301
$ws.Cells.Item(13, 7).Value = 4.603447666666667
$ws.Cells.Item(13, 8).Value = 13.810343
$ws.Cells.Item(13, 9).Value = 0.4142501574471451
$ws.Cells.Item(13, 10).Value = 0.4142501574471449
$ws.Cells.Item(13, 13).Value = 4.873057999999999
$ws.Cells.Item(13, 14).Value = 14.619174
$ws.Cells.Item(13, 15).Value = 0.07136450571781097
$ws.Cells.Item(13, 16).Value = 0.07136450571781099
$ws.Cells.Item(13, 17).Value = 22.43286747963133
$ws.Cells.Item(13, 18).Value = 201.895807316682
$ws.Cells.Item(13, 19).Value = 0.02956275772974088
$ws.Cells.Item(13, 20).Value = 0.02956275772974088
$ws.Cells.Item(14, 7).Value = 0.141604
$ws.Cells.Item(14, 8).Value = 0.424812
$ws.Cells.Item(14, 9).Value = 0.01274251029720526
$ws.Cells.Item(14, 10).Value = 0.01274251029720526
$ws.Cells.Item(14, 13).Value = 14.25737566666667
$ws.Cells.Item(14, 14).Value = 42.772127
$ws.Cells.Item(14, 15).Value = 0.2087950866344732
$ws.Cells.Item(14, 16).Value = 0.2087950866344732
$ws.Cells.Item(14, 17).Value = 2.018901423902667
$ws.Cells.Item(14, 18).Value = 18.170112815124
$ws.Cells.Item(14, 19).Value = 0.002660573541445639
$ws.Cells.Item(14, 20).Value = 0.002660573541445639
$ws.Cells.Item(15, 7).Value = 0.141604
$ws.Cells.Item(15, 8).Value = 0.424812
$ws.Cells.Item(15, 9).Value = 0.01274251029720526
$ws.Cells.Item(15, 10).Value = 0.01274251029720526
$ws.Cells.Item(15, 14).Value = 87.128332
$ws.Cells.Item(15, 15).Value = 0.4253229592313036
$ws.Cells.Item(15, 16).Value = 0.4253229592313036
$ws.Cells.Item(15, 17).Value = 4.112573441509333
$ws.Cells.Item(15, 18).Value = 37.013160973584
$ws.Cells.Item(15, 19).Value = 0.005419682187642699
$ws.Cells.Item(15, 20).Value = 0.005419682187642699
$ws.Cells.Item(16, 7).Value = 0.141604
$ws.Cells.Item(16, 8).Value = 0.424812
$ws.Cells.Item(16, 9).Value = 0.01274251029720526
$ws.Cells.Item(16, 10).Value = 0.01274251029720526
$ws.Cells.Item(16, 13).Value = 20.11084633333333
$ws.Cells.Item(16, 14).Value = 60.332539
$ws.Cells.Item(16, 15).Value = 0.2945174484164121
$ws.Cells.Item(16, 16).Value = 0.2945174484164122
$ws.Cells.Item(16, 17).Value = 2.847776284185333
$ws.Cells.Item(16, 18).Value = 25.629986557668
$ws.Cells.Item(16, 19).Value = 0.003752891619152751
$ws.Cells.Item(16, 20).Value = 0.003752891619152751
$ws.Cells.Item(17, 7).Value = 0.141604
$ws.Cells.Item(17, 8).Value = 0.424812
$ws.Cells.Item(17, 9).Value = 0.01274251029720526
$ws.Cells.Item(17, 10).Value = 0.01274251029720526
$ws.Cells.Item(17, 13).Value = 4.873057999999999
$ws.Cells.Item(17, 14).Value = 14.619174
$ws.Cells.Item(17, 15).Value = 0.07136450571781097
$ws.Cells.Item(17, 16).Value = 0.07136450571781099
$ws.Cells.Item(17, 17).Value = 0.690044505032
$ws.Cells.Item(17, 18).Value = 6.210400545288
$ws.Cells.Item(17, 19).Value = 0.00090936294896417
$ws.Cells.Item(17, 20).Value = 0.0009093629489641701
